$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last four data rows (6-9), shifting nothing up since they're at the end
$ws.Rows("6:9").Delete()

# Update the remaining data rows (2-5) with the new control point counts
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 1593

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 1458

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 1073

$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 134
